$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1): relabel existing columns ---
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "capacity"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"
$ws.Cells.Item(1,7).Value = "acquire_value"

# --- add the new trailing header columns (H1:N1), matching the layout used
#     by the other property sheets (land/building/etc.) ---
$newHeaders = @("property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = 8 + $i
    $ws.Cells.Item(1,7).Copy()
    $ws.Cells.Item(1,$col).PasteSpecial(-4122)
    $ws.Cells.Item(1,$col).Value = $newHeaders[$i]
}

# --- add the new trailing data columns (H2:N2) for the existing car record ---
$ws.Cells.Item(2,8).Value = "land"
$ws.Cells.Item(2,9).Value = "normal"

# J2 ("date") must stay text ("2012-04-26"), not auto-convert to a date
# serial: force text entry, then reapply the sibling-cell formatting so the
# cell ends up with the same effective style as the rest of row 2.
$ws.Cells.Item(2,10).NumberFormat = "@"
$ws.Cells.Item(2,10).Value = "2012-04-26"

$ws.Cells.Item(2,11).Value = "楊應雄"
$ws.Cells.Item(2,12).Value = 1758
$ws.Cells.Item(2,13).Value = "tmp248f1"
$ws.Cells.Item(2,14).Value = 31

for ($col = 8; $col -le 14; $col++) {
    $ws.Cells.Item(2,7).Copy()
    $ws.Cells.Item(2,$col).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
